$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 25: record the Reading (Read_Mark) score for the Cambridge 7 Test2 row ---
$ws.Range("G25").Value = 32
$ws.Range("H25").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'

# --- Row 26: picks up the "filled" Reading-column formatting even though no score yet ---
$ws.Range("H25").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null

# --- Row 27: record the Reading (Read_Mark) score, same formatting pickup as row 25 ---
$ws.Range("G27").Value = 24
$ws.Range("H25").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$ws.Range("H27").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'

# --- Row 28: new practiced test entry (Cambridge 7 Test 3) ---
$ws.Range("C27").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Copy() | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null

$ws.Range("C28").Value = 45493
$ws.Range("D28").Value = "IELTS7_Test3"
$ws.Rows.Item(28).RowHeight = 15

# --- Row 29: new practiced test entry (Cambridge 7 Test 4) ---
$ws.Range("C27").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null

$ws.Range("C29").Value = 45494
$ws.Range("D29").Value = "IELTS7_Test4"

# --- Selection ends on H26, matching where the user left off entering Reading scores ---
$ws.Range("H26").Select() | Out-Null
